# Auto-generated Excel COM-interop script
# Applies numeric cell updates across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# (Ixion Profits workbook) per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 550
$ws.Range("I12").Value = 550
$ws.Range("K12").Value = 550
$ws.Range("M12").Value = -380
$ws.Range("H69").Value = 5642.636
$ws.Range("I69").Value = 4766.6665
$ws.Range("J69").Value = 5971.125
$ws.Range("K69").Value = 14299.9995
$ws.Range("L69").Value = 17913.375
$ws.Range("M69").Value = -13425.9995
$ws.Range("N69").Value = -19661.375
$ws.Range("H72").Value = 5642.636
$ws.Range("I72").Value = 4766.6665
$ws.Range("J72").Value = 5971.125
$ws.Range("K72").Value = 42899.9985
$ws.Range("L72").Value = 53740.125
$ws.Range("M72").Value = -38531.9985
$ws.Range("N72").Value = -62476.125
$ws.Range("H137").Value = 1272.9193
$ws.Range("I137").Value = 758.75757
$ws.Range("J137").Value = 1858
$ws.Range("K137").Value = 2276.27271
$ws.Range("L137").Value = 5574
$ws.Range("M137").Value = 273.7272899999998
$ws.Range("N137").Value = -10674
$ws.Range("H138").Value = 2520.1064
$ws.Range("I138").Value = 996.65
$ws.Range("J138").Value = 3648.5925
$ws.Range("K138").Value = 2989.95
$ws.Range("L138").Value = 10945.7775
$ws.Range("M138").Value = 2150.05
$ws.Range("N138").Value = -21225.7775

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2082.077
$ws.Range("I2").Value = 2788.875
$ws.Range("J2").Value = 951.2
$ws.Range("K2").Value = 2788.875
$ws.Range("L2").Value = 951.2
$ws.Range("M2").Value = -2675.875
$ws.Range("N2").Value = -1177.2
$ws.Range("H63").Value = 166668530
$ws.Range("I63").Value = 200001630
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 200001630
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -200000944
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 166668530
$ws.Range("I66").Value = 200001630
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 1000008150
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -1000004718
$ws.Range("N66").Value = -21864
$ws.Range("H102").Value = 2850926.5
$ws.Range("I102").Value = 3368595
$ws.Range("K102").Value = 3368595
$ws.Range("M102").Value = -3366973
$ws.Range("H116").Value = 2082.077
$ws.Range("I116").Value = 2788.875
$ws.Range("J116").Value = 951.2
$ws.Range("K116").Value = 2788.875
$ws.Range("L116").Value = 951.2
$ws.Range("M116").Value = -494.875
$ws.Range("N116").Value = -5539.2
$ws.Range("H132").Value = 2600.125
$ws.Range("I132").Value = 1643.0968
$ws.Range("J132").Value = 3786.84
$ws.Range("K132").Value = 4929.2904
$ws.Range("L132").Value = 11360.52
$ws.Range("M132").Value = -2399.2904
$ws.Range("N132").Value = -16420.52

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2082.077
$ws.Range("I3").Value = 2788.875
$ws.Range("J3").Value = 951.2
$ws.Range("K3").Value = 2788.875
$ws.Range("L3").Value = 951.2
$ws.Range("M3").Value = -2674.875
$ws.Range("N3").Value = -1179.2
$ws.Range("H8").Value = 1235
$ws.Range("I8").Value = 1235
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1235
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H80").Value = 1700.875
$ws.Range("I80").Value = 1625
$ws.Range("J80").Value = 1776.75
$ws.Range("K80").Value = 1625
$ws.Range("L80").Value = 1776.75
$ws.Range("M80").Value = -627
$ws.Range("N80").Value = -3772.75
$ws.Range("H83").Value = 1700.875
$ws.Range("I83").Value = 1625
$ws.Range("J83").Value = 1776.75
$ws.Range("K83").Value = 8125
$ws.Range("L83").Value = 8883.75
$ws.Range("M83").Value = -3133
$ws.Range("N83").Value = -18867.75
$ws.Range("H134").Value = 1883.7805
$ws.Range("I134").Value = 1622.625
$ws.Range("J134").Value = 2252.4707
$ws.Range("K134").Value = 4867.875
$ws.Range("L134").Value = 6757.4121
$ws.Range("M134").Value = -2332.875
$ws.Range("N134").Value = -11827.4121

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2641.7747
$ws.Range("I31").Value = 2038.6842
$ws.Range("J31").Value = 2862.1345
$ws.Range("K31").Value = 2038.6842
$ws.Range("L31").Value = 2862.1345
$ws.Range("M31").Value = -1743.6842
$ws.Range("N31").Value = -3452.1345
$ws.Range("H34").Value = 2641.7747
$ws.Range("I34").Value = 2038.6842
$ws.Range("J34").Value = 2862.1345
$ws.Range("K34").Value = 2038.6842
$ws.Range("L34").Value = 2862.1345
$ws.Range("M34").Value = -1836.6842
$ws.Range("N34").Value = -3266.1345
$ws.Range("H86").Value = 125002660
$ws.Range("I86").Value = 166669550
$ws.Range("J86").Value = 1970
$ws.Range("K86").Value = 166669550
$ws.Range("L86").Value = 1970
$ws.Range("M86").Value = -166668427
$ws.Range("N86").Value = -4216
$ws.Range("H89").Value = 125002660
$ws.Range("I89").Value = 166669550
$ws.Range("J89").Value = 1970
$ws.Range("K89").Value = 833347750
$ws.Range("L89").Value = 9850
$ws.Range("M89").Value = -833342134
$ws.Range("N89").Value = -21082
$ws.Range("H99").Value = 12519280
$ws.Range("I99").Value = 20200
$ws.Range("K99").Value = 20200
$ws.Range("M99").Value = -18702
$ws.Range("H125").Value = 49750
$ws.Range("J125").Value = 49750
$ws.Range("L125").Value = 49750
$ws.Range("N125").Value = -54670
$ws.Range("H126").Value = 12519280
$ws.Range("I126").Value = 20200
$ws.Range("K126").Value = 60600
$ws.Range("M126").Value = -58130
$ws.Range("H132").Value = 2572.5
$ws.Range("I132").Value = 1637.409
$ws.Range("J132").Value = 6001.1665
$ws.Range("K132").Value = 4912.227000000001
$ws.Range("L132").Value = 18003.4995
$ws.Range("M132").Value = -2382.227000000001
$ws.Range("N132").Value = -23063.4995
$ws.Range("H135").Value = 33242.94
$ws.Range("J135").Value = 33242.94
$ws.Range("L135").Value = 33242.94
$ws.Range("N135").Value = -43382.94

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 40876.25
$ws.Range("J48").Value = 40876.25
$ws.Range("L48").Value = 122628.75
$ws.Range("N48").Value = -123128.75
$ws.Range("H58").Value = 3642.8572
$ws.Range("I58").Value = 500
$ws.Range("J58").Value = 4900
$ws.Range("K58").Value = 1500
$ws.Range("L58").Value = 14700
$ws.Range("M58").Value = -1372
$ws.Range("N58").Value = -14956
$ws.Range("H68").Value = 3535.434
$ws.Range("I68").Value = 4652
$ws.Range("J68").Value = 2375.923
$ws.Range("K68").Value = 13956
$ws.Range("L68").Value = 7127.768999999999
$ws.Range("M68").Value = -13145
$ws.Range("N68").Value = -8749.769
$ws.Range("H71").Value = 3535.434
$ws.Range("I71").Value = 4652
$ws.Range("J71").Value = 2375.923
$ws.Range("K71").Value = 41868
$ws.Range("L71").Value = 21383.307
$ws.Range("M71").Value = -37812
$ws.Range("N71").Value = -29495.307
$ws.Range("H107").Value = 764.67645
$ws.Range("I107").Value = 268.5
$ws.Range("K107").Value = 805.5
$ws.Range("M107").Value = 1114.5
$ws.Range("H122").Value = 649.8182
$ws.Range("I122").Value = 568
$ws.Range("J122").Value = 748
$ws.Range("K122").Value = 5112
$ws.Range("L122").Value = 6732
$ws.Range("M122").Value = -2662
$ws.Range("N122").Value = -11632
$ws.Range("H129").Value = 1715.1724
$ws.Range("I129").Value = 929.2143
$ws.Range("J129").Value = 2448.7334
$ws.Range("K129").Value = 2787.6429
$ws.Range("L129").Value = 7346.2002
$ws.Range("M129").Value = 2212.3571
$ws.Range("N129").Value = -17346.2002
$ws.Range("H131").Value = 16177664
$ws.Range("J131").Value = 18869338
$ws.Range("L131").Value = 56608014
$ws.Range("N131").Value = -56618094
$ws.Range("H132").Value = 1854368.8
$ws.Range("I132").Value = 841
$ws.Range("J132").Value = 2139527
$ws.Range("K132").Value = 7569
$ws.Range("L132").Value = 19255743
$ws.Range("M132").Value = -5039
$ws.Range("N132").Value = -19260803
$ws.Range("H133").Value = 33586.95
$ws.Range("I133").Value = 115676.78
$ws.Range("J133").Value = 8960
$ws.Range("K133").Value = 347030.34
$ws.Range("L133").Value = 26880
$ws.Range("M133").Value = -341970.34
$ws.Range("N133").Value = -37000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6110.1816
$ws.Range("I126").Value = 6358.2856
$ws.Range("J126").Value = 900
$ws.Range("K126").Value = 19074.8568
$ws.Range("L126").Value = 2700
$ws.Range("M126").Value = -16604.8568
$ws.Range("N126").Value = -7640
$ws.Range("H132").Value = 3297.3142
$ws.Range("I132").Value = 2629.7917
$ws.Range("J132").Value = 3645.587
$ws.Range("K132").Value = 7889.375100000001
$ws.Range("L132").Value = 10936.761
$ws.Range("M132").Value = -5359.375100000001
$ws.Range("N132").Value = -15996.761

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2682.8235
$ws.Range("I7").Value = 2178.6667
$ws.Range("J7").Value = 3250
$ws.Range("K7").Value = 2178.6667
$ws.Range("L7").Value = 3250
$ws.Range("M7").Value = -2066.6667
$ws.Range("N7").Value = -3474
$ws.Range("H93").Value = 55578944
$ws.Range("I93").Value = 29357.572
$ws.Range("K93").Value = 29357.572
$ws.Range("M93").Value = -28109.572
$ws.Range("H126").Value = 2682.8235
$ws.Range("I126").Value = 2178.6667
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 6536.000100000001
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = -4066.000100000001
$ws.Range("N126").Value = -14690
$ws.Range("H132").Value = 9632281
$ws.Range("I132").Value = 12747445
$ws.Range("J132").Value = 3591.182
$ws.Range("K132").Value = 38242335
$ws.Range("L132").Value = 10773.546
$ws.Range("M132").Value = -38239805
$ws.Range("N132").Value = -15833.546
$ws.Range("H135").Value = 25700
$ws.Range("J135").Value = 22600
$ws.Range("L135").Value = 22600
$ws.Range("N135").Value = -32740

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 14021.8
$ws.Range("I122").Value = 15701.333
$ws.Range("J122").Value = 11502.5
$ws.Range("K122").Value = 47103.999
$ws.Range("L122").Value = 34507.5
$ws.Range("M122").Value = -44653.999
$ws.Range("N122").Value = -39407.5
$ws.Range("H136").Value = 2450.7073
$ws.Range("I136").Value = 2673.4614
$ws.Range("J136").Value = 2064.6
$ws.Range("K136").Value = 8020.3842
$ws.Range("L136").Value = 6193.799999999999
$ws.Range("M136").Value = -5470.3842
$ws.Range("N136").Value = -11293.8
